# Apply price / ranking updates per the "Updated symbol list" commit.
# Column D holds numeric-looking values that are stored as TEXT in the
# workbook, so those assignments are prefixed with a leading apostrophe
# to force text entry (matching the original inline-string cell type)
# instead of letting Excel coerce them into Number cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'268.99"
$ws.Range("D4").Value = "'6.250"
$ws.Range("D5").Value = "'0.06207"
$ws.Range("D7").Value = "'6.537"
$ws.Range("D8").Value = "'1.388"
$ws.Range("D9").Value = "'0.8270"
$ws.Range("D10").Value = "'0.1639"
$ws.Range("D11").Value = "'0.08297"
$ws.Range("D13").Value = "'0.03192"
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D14").Value = "'0.1242"
$ws.Range("E14").Value = "13ProBitTokenPROB"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09200"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.770"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001631"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04691"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006338"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.006204"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.001068"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.725"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.291"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.01366"
$ws.Range("E25").Value = "24OneONE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3289"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("D40").Value = "'0.04732"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1122"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003521"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.01164"
$ws.Range("D45").Value = "'0.00006293"
$ws.Range("D46").Value = "'0.0009902"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D49").Value = "'0.002338"
$ws.Range("E49").Value = "48BOLOBOLO"
